$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.691.59'
$ws.Range("E2").Value = '  -3.69%  '

$ws.Range("D3").Value = '3.821.63'
$ws.Range("E3").Value = '  -3.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.03%  '

$ws.Range("D7").Value = '3.821.60'
$ws.Range("E7").Value = '  -3.06%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -2.17%  '

$ws.Range("E10").Value = '  -5.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000258'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.89%  '

$ws.Range("D15").Value = '4.464.88'
$ws.Range("E15").Value = '  -3.15%  '

$ws.Range("D16").Value = '3.817.60'
$ws.Range("E16").Value = '  -3.85%  '

$ws.Range("D17").Value = '67.912.09'
$ws.Range("E17").Value = '  -3.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.82%  '

$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.42%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '463.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.19%  '

$ws.Range("E23").Value = '  -2.70%  '

$ws.Range("E24").Value = '  -5.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.92%  '

$ws.Range("E27").Value = '  -3.34%  '

$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.57%  '

$ws.Range("D31").Value = '3.970.86'
$ws.Range("E31").Value = '  -3.09%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.49%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.62%  '

$ws.Range("D36").Value = '3.784.04'
$ws.Range("E36").Value = '  -3.20%  '

$ws.Range("E37").Value = '  -4.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.25%  '

$ws.Range("E39").Value = '  -1.30%  '

$ws.Range("E40").Value = '  -2.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.37%  '

$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.310'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.12%  '

$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '417.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.80%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.03%  '

$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000292'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0357'
$ws.Range("D51").Style = "Normal"
